# Populate the navigation-components worksheet (Sheet1) with the full
# component/dropdown-source table that replaces the original 2x2 sample.
#
# The fill order below intentionally matches how the shared-string table
# was built in the target workbook: the original A1:B2 four cells stay as
# they were, then the small "dropdown source" list in column D (rows 1-3)
# is written, then row 3 of the main table is completed, and finally the
# remaining rows of the main table are filled in top-to-bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing header / first data row (values unchanged, already present) ---
$ws.Range("A1").Value = "ComponentName"
$ws.Range("B1").Value = "Type"
$ws.Range("A2").Value = "Pclass"
$ws.Range("B2").Value = "C"

# --- dropdown-source helper list in column D (rows 1-3) ---
$ws.Range("D1").Value = "C- Click"
$ws.Range("D2").Value = "D- Dropdown"

# --- finish row 3 of the main table ---
$ws.Range("A3").Value = "Cclass"
$ws.Range("B3").Value = "C"
$ws.Range("D3").Value = "NA - Not Applicable"

# --- remaining main-table rows (4-14), column A name / column B type ---
$ws.Range("A4").Value = "Cuvette"
$ws.Range("B4").Value = "C"
$ws.Range("A5").Value = "NanoVolume"
$ws.Range("B5").Value = "C"
$ws.Range("A6").Value = "Cuvette"
$ws.Range("B6").Value = "C"
$ws.Range("A7").Value = "StoreMethods"
$ws.Range("B7").Value = "C"
$ws.Range("A8").Value = "DataFiles"
$ws.Range("B8").Value = "C"
$ws.Range("A9").Value = "NucleicAcid"
$ws.Range("B9").Value = "C"
$ws.Range("A10").Value = "DyeLabeledNucleicAcid"
$ws.Range("B10").Value = "C"
$ws.Range("A11").Value = "ProteinUV"
$ws.Range("B11").Value = "C"
$ws.Range("A12").Value = "DyeLabeledProtein"
$ws.Range("B12").Value = "C"
$ws.Range("A13").Value = "GeneralMethods"
$ws.Range("B13").Value = "C"
$ws.Range("A14").Value = "NucleicAcid_SubTypes"
$ws.Range("B14").Value = "C"

# --- apply a uniform style touch across the whole used range A1:D14 ---
# (this both stamps every cell - including the blank column C spacer and
# the blank column D cells in rows 4-14 - with a style index, and forces
# column C cells to exist in the sheet even though they hold no value)
$ws.Range("A1:D14").WrapText = $false

# --- final selection, matching the saved workbook view state ---
$null = $ws.Range("F14").Select()
